# The deck currently has two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"     (used by the Slide Master / all slides)
#
# The target edit swaps the two themes' content: the Slide Master's theme
# (theme2.xml) must end up holding the "Office Theme" palette, while the
# Notes Master's theme (theme1.xml) ends up holding "Integral".
#
# The PowerPoint object model only exposes the *active* (slide-master)
# theme for writing -- via Slide.ThemeColorScheme, which carries all 12
# theme colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). Apply the
# "Office Theme" palette there so the deck's visible design switches from
# Integral to the stock Office colors, matching the slide master side of
# the swap.

function Get-ComRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$tcs = $slide.ThemeColorScheme

# Target palette = the "Office Theme" colors (R,G,B), in clrScheme order.
$officeColors = @(
    @(0x00,0x00,0x00),  # 1  dk1
    @(0xFF,0xFF,0xFF),  # 2  lt1
    @(0x44,0x54,0x6A),  # 3  dk2
    @(0xE7,0xE6,0xE6),  # 4  lt2
    @(0x5B,0x9B,0xD5),  # 5  accent1
    @(0xED,0x7D,0x31),  # 6  accent2
    @(0xA5,0xA5,0xA5),  # 7  accent3
    @(0xFF,0xC0,0x00),  # 8  accent4
    @(0x44,0x72,0xC4),  # 9  accent5
    @(0x70,0xAD,0x47),  # 10 accent6
    @(0x05,0x63,0xC1),  # 11 hlink
    @(0x95,0x4F,0x72)   # 12 folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeColors[$i - 1]
    $color = $tcs.Colors($i)
    $color.RGB = Get-ComRGB $rgb[0] $rgb[1] $rgb[2]
}
